$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the B1:B6 values (cells removed entirely from the sheet)
$ws.Range("B1:B6").ClearContents()

# Update the report generation date/time in B11
$ws.Range("B11").Value = "07-07-2022 12:12:53 pm"

# Insert a new row above row 15 ("Fin" row) to host the new "Componente" row,
# which pushes the existing "Fin" row down to row 16
$ws.Rows.Item(15).Insert()
$ws.Range("A15:H15").Style = "Normal"

# Populate the newly inserted row 15 with the "Componente" data.
$ws.Range("A15").Value = "Componente"
$ws.Range("B15").Value = "6 Gobierno Ciudadano"

# C15 and F15 hold numeric-looking text ("5612"/"70"). Assigning them directly via
# .Value would be auto-coerced to numbers, so instead build the text in a scratch
# cell via TEXT() and copy/paste the value, which preserves the text type to match
# the original workbook's inline-string cells.
$ws.Range("Z1").Formula = '=TEXT(5612,"0")'
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("C15").PasteSpecial(-4163) | Out-Null

$ws.Range("D15").Value = "Sumatoria de POA alineados con el Plan Estatal de Desarrollo 2021-2027"

# Row 15 has no value in column E (unlike row 16 below it) - drop the placeholder
# cell that Insert() left behind.
$ws.Range("E15").ClearContents()

$ws.Range("Z1").Formula = '=TEXT(70,"0")'
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("F15").PasteSpecial(-4163) | Out-Null

$ws.Range("G15").Value = "Mensual"
$ws.Range("H15").Value = "Plan"

# Clean up the scratch cell used for text coercion
$ws.Range("Z1").ClearContents()
